$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.072.78"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.23%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.957.56"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "380.31"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.27%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.19"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.48%  "

$ws.Range("E7").Value = "  +1.54%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.591"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.50"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.04%  "

$ws.Range("E11").Value = "  -1.30%  "

$ws.Range("E12").Value = "  +1.78%  "

$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.422.86"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.41%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.80"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.57%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "18.33"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "11.52"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +46.85%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.962.21"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.12%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.997"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.15%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "51.159.19"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("E20").Value = "  -1.06%  "

$ws.Range("E21").Value = "  -2.43%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0959"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.31"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +14.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.12"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "266.64"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.71%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.87"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -7.70%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.18"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -10.56%  "

$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.166"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.20%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "25.84"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.53%  "

$ws.Range("E31").Value = "  -2.67%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "10.32"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +4.22%  "

$ws.Range("B33").Value = "Toncoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.12"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.64%  "

$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "34.34"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.37%  "

$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "51.08"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.34%  "

$ws.Range("E36").Value = "  -3.09%  "

$ws.Range("E37").Value = "  +0.14%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.23"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +7.65%  "

$ws.Range("E39").Value = "  +0.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.53"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.19%  "

$ws.Range("E41").Value = "  +1.63%  "

$ws.Range("B42").Value = "Monero"
$ws.Range("C42").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "125.36"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.68%  "

$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.50"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.92%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.50"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.53"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +7.83%  "

$ws.Range("E46").Value = "  -0.76%  "

$ws.Range("E47").Value = "  +3.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.270"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.45%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.051.75"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.56%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0322"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -7.31%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "5.41"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.82%  "
